$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (pushes old Visit/Specimen_Number columns right)
$ws.Columns("E:E").Insert()

# New header cell: "Material Type" with a bold white-on-black style
$ws.Range("E1").Value = "Material Type"
$ws.Range("E1").Interior.Color = 0
$ws.Range("E1").Font.Color = 16777215
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108

# Fill the new column's data rows with the constant "DNA:Genomic"
for ($r = 2; $r -le 24; $r++) {
  $ws.Cells.Item($r, 5).Value = "DNA:Genomic"
  $ws.Cells.Item($r, 5).HorizontalAlignment = -4108
}

# Widen the new column to match the neighboring Collection_Date column
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Restore selection to reflect the newly inserted column
$ws.Range("E1:E24").Select() | Out-Null
